$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "site_url" column (column I),
# shifting site_url and its data to column J.
$ws.Columns.Item(9).Insert() | Out-Null

# Populate the new "node_id" optional column.
$ws.Range("I1").Value = "node_id"
$ws.Range("I2").Value = "img"

# Match the author's final selection in the worksheet.
$ws.Range("I7").Select() | Out-Null
